# Recognize individual words in the "Abkürzungen" (abbreviations) column:
#  - A6 "DFA,DA"  -> "DFA, DA"   (space added after comma)
#  - A7 "Berlin"  -> "Berlin, TXL" (second abbreviation appended)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "DFA, DA"
$ws.Range("A7").Value = "Berlin, TXL"

# Column A grows a bit wider to fit the longer text; select the cell that
# was being edited (A7) as the last user action.
$ws.Columns("A").ColumnWidth = 35.16666666666667
$ws.Range("A7").Select()
